$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8308
$ws1.Range("F5").Value = 6058
$ws1.Range("F7").Value = 105
$ws1.Range("F8").Value = 19
$ws1.Range("F10").Value = 313
$ws1.Range("F11").Value = 985

# Sheet "全部类型" (sheet4): update 想去人数 (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8309
$ws4.Range("F5").Value = 6058
$ws4.Range("F7").Value = 105
$ws4.Range("F8").Value = 19
$ws4.Range("F10").Value = 313
$ws4.Range("F15").Value = 985
